# edit.ps1 - applies t10.1 sheet update:
#  1) revises D-column values for existing rows 2-131 (2010-2022 data revision)
#  2) appends a new 2023 year block as rows 132-141

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Revised "Valor" (column D) figures for existing rows ---
$ws.Range("D2").Value = 7948.173962135745
$ws.Range("D3").Value = 3.864969896044429
$ws.Range("D5").Value = 2228.978462293759
$ws.Range("D6").Value = 1861.941910616449
$ws.Range("D7").Value = 987.5962546332948
$ws.Range("D8").Value = 78.57136469340311
$ws.Range("D9").Value = 2227.558485104306
$ws.Range("D10").Value = 439.7287869325916
$ws.Range("D11").Value = 123.7986978619414
$ws.Range("D12").Value = 8260.960460759228
$ws.Range("D14").Value = 0.3761258600430559
$ws.Range("D15").Value = 2201.139126414258
$ws.Range("D16").Value = 2093.572101459282
$ws.Range("D17").Value = 1036.732797375233
$ws.Range("D18").Value = 86.14913781057841
$ws.Range("D19").Value = 2223.657458749666
$ws.Range("D20").Value = 447.6921470034345
$ws.Range("D21").Value = 172.0176919467761
$ws.Range("D22").Value = 8905.919549363491
$ws.Range("D23").Value = 3.742007350094518
$ws.Range("D25").Value = 2299.613213787246
$ws.Range("D26").Value = 2446.752269782316
$ws.Range("D27").Value = 1163.528153964236
$ws.Range("D28").Value = 88.31577267755924
$ws.Range("D29").Value = 2246.169122874369
$ws.Range("D30").Value = 455.6520039500786
$ws.Range("D31").Value = 205.8890123276872
$ws.Range("D32").Value = 9485.294132059897
$ws.Range("D33").Value = 3.748027280592707
$ws.Range("D35").Value = 2470.038068519863
$ws.Range("D36").Value = 2641.53784223107
$ws.Range("D37").Value = 1361.059100962156
$ws.Range("D38").Value = 113.6852163291302
$ws.Range("D39").Value = 2272.795398729881
$ws.Range("D40").Value = 409.216361693038
$ws.Range("D41").Value = 216.962143594757
$ws.Range("D42").Value = 10178.74234699472
$ws.Range("D43").Value = 3.737023563155108
$ws.Range("D45").Value = 2368.517352943205
$ws.Range("D46").Value = 3228.43698794786
$ws.Range("D47").Value = 1421.509814907411
$ws.Range("D48").Value = 126.6094611710749
$ws.Range("D49").Value = 2229.627162764234
$ws.Range("D50").Value = 537.9212894754992
$ws.Range("D51").Value = 266.1202777854409
$ws.Range("D52").Value = 9356.577031207362
$ws.Range("D53").Value = 3.777523282056657
$ws.Range("D54").Value = 0.3923570460908749
$ws.Range("D55").Value = 2245.031701925366
$ws.Range("D56").Value = 2810.122342165998
$ws.Range("D57").Value = 1294.027009955109
$ws.Range("D58").Value = 107.0280885722107
$ws.Range("D59").Value = 2117.245901559477
$ws.Range("D60").Value = 461.2659022538651
$ws.Range("D61").Value = 321.8560847753367
$ws.Range("D62").Value = 9131.98252761689
$ws.Range("D65").Value = 2130.33573919216
$ws.Range("D66").Value = 2641.559698942563
$ws.Range("D67").Value = 1358.530153634022
$ws.Range("D68").Value = 85.96040819775109
$ws.Range("D69").Value = 2141.540134987541
$ws.Range("D70").Value = 462.9032300353873
$ws.Range("D71").Value = 311.1531626274681
$ws.Range("D72").Value = 9038.521203142096
$ws.Range("D75").Value = 2113.820236173995
$ws.Range("D76").Value = 2570.336699922301
$ws.Range("D77").Value = 1371.130958868779
$ws.Range("D78").Value = 94.45595768906894
$ws.Range("D79").Value = 2117.874725300317
$ws.Range("D80").Value = 449.9174141679175
$ws.Range("D81").Value = 320.9852110197174
$ws.Range("D82").Value = 8905.8339818159
$ws.Range("D83").Value = 3.618636503777034
$ws.Range("D84").Value = 0.3720399211117291
$ws.Range("D85").Value = 2047.894952974252
$ws.Range("D86").Value = 2761.332443897886
$ws.Range("D87").Value = 1433.587684808288
$ws.Range("D88").Value = 93.12796412585945
$ws.Range("D89").Value = 1809.986402155821
$ws.Range("D90").Value = 486.7353694401258
$ws.Range("D91").Value = 273.1691644136682
$ws.Range("D92").Value = 9036.115716607301
$ws.Range("D94").Value = 0.3633506477362351
$ws.Range("D95").Value = 2440.652249760819
$ws.Range("D96").Value = 2595.178671036135
$ws.Range("D97").Value = 1390.091050975048
$ws.Range("D98").Value = 102.586029477938
$ws.Range("D99").Value = 1663.888897455947
$ws.Range("D100").Value = 516.7889475156243
$ws.Range("D101").Value = 326.9298703857901
$ws.Range("D102").Value = 7762.529020645291
$ws.Range("D103").Value = 3.434674100231558
$ws.Range("D105").Value = 1750.654577325427
$ws.Range("D106").Value = 2551.742644372269
$ws.Range("D107").Value = 956.2732850309048
$ws.Range("D108").Value = 131.4423096164256
$ws.Range("D109").Value = 1437.455748517693
$ws.Range("D110").Value = 643.6433018114983
$ws.Range("D111").Value = 291.3171539710714
$ws.Range("D112").Value = 8494.877891523587
$ws.Range("D113").Value = 3.384235570539732
$ws.Range("D114").Value = 0.328446387139379
$ws.Range("D115").Value = 2254.498066292051
$ws.Range("D116").Value = 2628.87594444971
$ws.Range("D117").Value = 1240.418899762093
$ws.Range("D118").Value = 141.5880861742734
$ws.Range("D119").Value = 1414.191982255018
$ws.Range("D120").Value = 510.4721256838192
$ws.Range("D121").Value = 304.8327869066216
$ws.Range("D122").Value = 10476.77956786359
$ws.Range("D123").Value = 3.452399617681198
$ws.Range("D124").Value = 0.3434005647028749
$ws.Range("D125").Value = 2632.70645710208
$ws.Range("D126").Value = 3302.587567009418
$ws.Range("D127").Value = 1728.694632294269
$ws.Range("D128").Value = 170.6266987792879
$ws.Range("D129").Value = 1424.726134429704
$ws.Range("D130").Value = 922.5716239147124
$ws.Range("D131").Value = 294.8664543341227

# --- 2) New rows for year 2023 ---
$ws.Range("A132").Value = 2023
$ws.Range("B132").Value = "Receita bruta de prestação de Serviços"
$ws.Range("C132").Value = "R`$ milhões"
$ws.Range("D132").Value = 11209.714

$ws.Range("A133").Value = 2023
$ws.Range("B133").Value = "Participação da receita bruta de Serviços em Sergipe no Nordeste"
$ws.Range("C133").Value = "%"
$ws.Range("D133").Value = 3.219876533226341

$ws.Range("A134").Value = 2023
$ws.Range("B134").Value = "Participação da receita bruta de Serviços em Sergipe no Brasil"
$ws.Range("C134").Value = "%"
$ws.Range("D134").Value = 0.3255402456283393

$ws.Range("A135").Value = 2023
$ws.Range("B135").Value = "   Transportes, serviços auxiliares aos transportes e correio"
$ws.Range("C135").Value = "R`$ milhões"
$ws.Range("D135").Value = 2752.555

$ws.Range("A136").Value = 2023
$ws.Range("B136").Value = "   Serviços profissionais, administrativos e complementares"
$ws.Range("C136").Value = "R`$ milhões"
$ws.Range("D136").Value = 3505.843

$ws.Range("A137").Value = 2023
$ws.Range("B137").Value = "   Serviços prestados principalmente às famílias"
$ws.Range("C137").Value = "R`$ milhões"
$ws.Range("D137").Value = 2044.616

$ws.Range("A138").Value = 2023
$ws.Range("B138").Value = "   Serviços de manutenção e reparação"
$ws.Range("C138").Value = "R`$ milhões"
$ws.Range("D138").Value = 212.594

$ws.Range("A139").Value = 2023
$ws.Range("B139").Value = "   Serviços de informação e comunicação"
$ws.Range("C139").Value = "R`$ milhões"
$ws.Range("D139").Value = 1546.602

$ws.Range("A140").Value = 2023
$ws.Range("B140").Value = "   Outras atividades de serviços"
$ws.Range("C140").Value = "R`$ milhões"
$ws.Range("D140").Value = 911.533

$ws.Range("A141").Value = 2023
$ws.Range("B141").Value = "   Atividades imobiliárias "
$ws.Range("C141").Value = "R`$ milhões"
$ws.Range("D141").Value = 235.971

